# Apply data refresh to 杭州-漫展信息.xlsx per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("C2").Value = "杭州·CICAF·中国国风品牌盛典、中国COSPLAY超级盛典"
$ws.Range("F2").Value = 658
$ws.Range("G2").Value = "不可售"
$ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=84828"
$ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202405/X6g42ZCh1715399314306.jpeg"
$ws.Range("B3").Value = "'2024-05-29"
$ws.Range("C3").Value = "杭州·第二十届中国国际动漫节主会场门票"
$ws.Range("D3").Value = "长江南路336号 白马湖国际会展中心"
$ws.Range("E3").Value = "2024.05.29 09:00-06.02 17:00"
$ws.Range("F3").Value = 6457
$ws.Range("G3").Value = 70
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84823"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202405/KfxshzO81715239999885.jpeg"
$ws.Range("F4").Value = 1032
$ws.Range("F5").Value = 637
$ws.Range("F6").Value = 1425
$ws.Range("F7").Value = 3175
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 527
$ws.Range("F10").Value = 2073
$ws.Range("F11").Value = 447
$ws.Range("F13").Value = 218
$ws.Range("F15").Value = 226
$ws.Range("F16").Value = 1029
$ws.Range("F17").Value = 397
$ws.Range("F18").Value = 63
$ws.Range("F19").Value = 147
$ws.Range("F20").Value = 3963
$ws.Range("F21").Value = 1211
$ws.Range("F22").Value = 3133
$ws.Range("F23").Value = 306
$ws.Range("F24").Value = 84
$ws.Range("F25").Value = 2889
$ws.Range("F26").Value = 2889
$ws.Range("F27").Value = 4529
$ws.Range("F30").Value = 502
$ws.Range("F31").Value = 2992
$ws.Range("F32").Value = 273
$ws.Range("F33").Value = 38
$ws.Range("F34").Value = 103
$ws.Range("F35").Value = 60
$ws.Range("F36").Value = 547
$ws.Range("F37").Value = 1085
$ws.Range("F38").Value = 1337
$ws.Range("F39").Value = 95
$ws.Range("F40").Value = 1200
$ws.Range("F41").Value = 775
$ws.Range("F43").Value = 710
$ws.Range("F44").Value = 467
$ws.Range("F45").Value = 39
$ws.Range("F46").Value = 176
$ws.Range("F47").Value = 23
$ws.Range("F48").Value = 71
$ws.Range("F49").Value = 339
$ws.Range("F50").Value = 3655

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F3").Value = 11
$ws.Range("F10").Value = 950

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1414

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6457
$ws.Range("C3").Value = "杭州·六一特献｜【直到世界尽头】灌篮高手等神级中日动漫演唱会，全体起立！"
$ws.Range("D3").Value = "萧山区杭州国际博览中心A座6楼 九莱福音乐现场"
$ws.Range("E3").Value = "2024.06.01 20:00-06.01 22:00"
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 128
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85097"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202405/xfdusgJP1715147982566.jpeg"
$ws.Range("B4").Value = "'2024-06-05"
$ws.Range("C4").Value = "杭州·英雄时代2024哈瓦西钢琴演奏会"
$ws.Range("D4").Value = "中国杭州北山路86号西湖岳湖景区 中国杭州西湖岳湖景区印象西湖"
$ws.Range("E4").Value = "2024.06.05 20:00-06.05 21:30"
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 499
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83902"
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/BFRFmKpT1712569969076.jpeg"
$ws.Range("F5").Value = 637
$ws.Range("F6").Value = 1425
$ws.Range("F7").Value = 3175
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 527
$ws.Range("F11").Value = 2073
$ws.Range("F12").Value = 447
$ws.Range("F14").Value = 218
$ws.Range("F15").Value = 950
$ws.Range("F18").Value = 226
$ws.Range("F19").Value = 1029
$ws.Range("F21").Value = 397
$ws.Range("F22").Value = 147
$ws.Range("F23").Value = 3964
$ws.Range("F25").Value = 1211
$ws.Range("F27").Value = 3133
$ws.Range("F28").Value = 2889
$ws.Range("F29").Value = 2889
$ws.Range("F30").Value = 4529
$ws.Range("F32").Value = 2992
$ws.Range("F33").Value = 273
$ws.Range("B34").Value = "'2024-07-20"
$ws.Range("C34").Value = "杭州·首届次元格子动漫展嘉宾内场——文森"
$ws.Range("D34").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws.Range("E34").Value = "2024.07.20 09:00-07.20 17:00"
$ws.Range("F34").Value = 103
$ws.Range("G34").Value = 238
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=86518"
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202405/kwjuTLK31716953514797.jpeg"
$ws.Range("C35").Value = "杭州·ACG CLUB动漫游戏嘉年华"
$ws.Range("D35").Value = "中心路1号 白蓝地文创街区"
$ws.Range("E35").Value = "2024.07.27 10:00-07.27 17:00"
$ws.Range("F35").Value = 547
$ws.Range("G35").Value = 68.88
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=86265"
$ws.Range("I35").Value = "//i1.hdslb.com/bfs/openplatform/202405/XBRfeQwu1716533419093.jpeg"
$ws.Range("B36").Value = "'2024-07-27"
$ws.Range("C36").Value = "杭州·夏之誓国乙only-日夜场"
$ws.Range("D36").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws.Range("E36").Value = "2024.07.27 10:00-07.27 21:00"
$ws.Range("F36").Value = 1085
$ws.Range("G36").Value = 69
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"
$ws.Range("C37").Value = "杭州·火影忍者only"
$ws.Range("D37").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E37").Value = "2024.07.28 09:00-07.28 18:00"
$ws.Range("F37").Value = 1337
$ws.Range("G37").Value = 75
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=84243"
$ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202405/isG309e51715657222196.jpeg"
$ws.Range("B38").Value = "'2024-07-28"
$ws.Range("C38").Value = "杭州·第三届百合Only·同好交流"
$ws.Range("D38").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E38").Value = "2024.07.28 10:00-07.28 16:00"
$ws.Range("F38").Value = 95
$ws.Range("G38").Value = 46
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=85895"
$ws.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202405/fP3O6LYz1715844791584.jpeg"
$ws.Range("C39").Value = "杭州·AP动漫游戏嘉年华"
$ws.Range("D39").Value = "沈半路171号 Tcar汽车文化主题公园"
$ws.Range("E39").Value = "2024.08.03 09:00-08.04 17:00"
$ws.Range("F39").Value = 1200
$ws.Range("G39").Value = 70
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=85527"
$ws.Range("I39").Value = "//i2.hdslb.com/bfs/openplatform/202405/JbVl16OE1715676665714.jpeg"
$ws.Range("B40").Value = "'2024-08-03"
$ws.Range("C40").Value = "杭州·梦漫星河动漫展"
$ws.Range("D40").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E40").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F40").Value = 775
$ws.Range("G40").Value = 68
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$ws.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"
$ws.Range("C41").Value = "杭州·【七夕巨献·早鸟6折】真的爱你”致敬Beyond·黄家驹31周年演唱会·630乐团再现91殿堂级演出"
$ws.Range("D41").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws.Range("E41").Value = "2024.08.10 19:30-08.10 21:30"
$ws.Range("F41").Value = 2
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=85333"
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202405/uYt32zt21715221330023.jpeg"
$ws.Range("C42").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D42").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E42").Value = "2024.08.10 10:00-08.10 17:00"
$ws.Range("F42").Value = 467
$ws.Range("G42").Value = 60
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("C43").Value = "杭州·吹响号角狩猎吧—怪物猎人&最终幻想&塞尔达 热血游戏视听音乐会"
$ws.Range("E43").Value = "2024.08.10 19:30-08.10 21:00"
$ws.Range("F43").Value = 44
$ws.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=84879"
$ws.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202404/jJcWHuGa1714010818763.jpeg"
$ws.Range("F44").Value = 39
$ws.Range("F46").Value = 176
$ws.Range("F47").Value = 23
$ws.Range("F48").Value = 71
$ws.Range("F49").Value = 339
$ws.Range("F50").Value = 3655
